$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"2"
$ws.Range("G2").Value = [double]"2.3842205"
$ws.Range("H2").Value = [double]"4.768441"
$ws.Range("I2").Value = [double]"0.0684902599354226"
$ws.Range("J2").Value = [double]"0.05735520746201143"
$ws.Range("K2").Value = [double]"2"
$ws.Range("M2").Value = [double]"35.789624"
$ws.Range("N2").Value = [double]"71.57924800000001"
$ws.Range("O2").Value = [double]"0.258139457682779"
$ws.Range("P2").Value = [double]"0.1993778771086309"
$ws.Range("Q2").Value = [double]"85.33035522809202"
$ws.Range("R2").Value = [double]"341.3214209123681"
$ws.Range("S2").Value = [double]"0.01768003855628255"
$ws.Range("T2").Value = [double]"0.01143535950490094"

$ws.Range("E3").Value = [double]"2"
$ws.Range("G3").Value = [double]"2.3842205"
$ws.Range("H3").Value = [double]"4.768441"
$ws.Range("I3").Value = [double]"0.0684902599354226"
$ws.Range("J3").Value = [double]"0.05735520746201143"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"60.113367"
$ws.Range("N3").Value = [double]"180.340101"
$ws.Range("O3").Value = [double]"0.4335790718803266"
$ws.Range("P3").Value = [double]"0.5023219368682956"
$ws.Range("Q3").Value = [double]"143.3235219254235"
$ws.Range("R3").Value = [double]"859.941131552541"
$ws.Range("S3").Value = [double]"0.02969594333564285"
$ws.Range("T3").Value = [double]"0.0288107789018005"

$ws.Range("E4").Value = [double]"2"
$ws.Range("G4").Value = [double]"2.3842205"
$ws.Range("H4").Value = [double]"4.768441"
$ws.Range("I4").Value = [double]"0.0684902599354226"
$ws.Range("J4").Value = [double]"0.05735520746201143"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"0.2072186666666667"
$ws.Range("N4").Value = [double]"0.621656"
$ws.Range("O4").Value = [double]"0.001494603973349423"
$ws.Range("P4").Value = [double]"0.001731569652308208"
$ws.Range("Q4").Value = [double]"0.4940549930493333"
$ws.Range("R4").Value = [double]"2.964329958296"
$ws.Range("S4").Value = [double]"0.0001023658146352174"
$ws.Range("T4").Value = [double]"9.931453664306029E-05"

$ws.Range("E5").Value = [double]"2"
$ws.Range("G5").Value = [double]"2.3842205"
$ws.Range("H5").Value = [double]"4.768441"
$ws.Range("I5").Value = [double]"0.0684902599354226"
$ws.Range("J5").Value = [double]"0.05735520746201143"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"21.40334366666667"
$ws.Range("N5").Value = [double]"64.210031"
$ws.Range("O5").Value = [double]"0.1543756795743782"
$ws.Range("P5").Value = [double]"0.178851553034748"
$ws.Range("Q5").Value = [double]"51.03029073861184"
$ws.Range("R5").Value = [double]"306.181744431671"
$ws.Range("S5").Value = [double]"0.01057323042175667"
$ws.Range("T5").Value = [double]"0.01025806792921091"

$ws.Range("E6").Value = [double]"2"
$ws.Range("G6").Value = [double]"2.3842205"
$ws.Range("H6").Value = [double]"4.768441"
$ws.Range("I6").Value = [double]"0.0684902599354226"
$ws.Range("J6").Value = [double]"0.05735520746201143"
$ws.Range("K6").Value = [double]"2"
$ws.Range("M6").Value = [double]"21.1309775"
$ws.Range("N6").Value = [double]"42.261955"
$ws.Range("O6").Value = [double]"0.1524111868891667"
$ws.Range("P6").Value = [double]"0.1177170633360173"
$ws.Range("Q6").Value = [double]"50.38090974053875"
$ws.Range("R6").Value = [double]"201.523638962155"
$ws.Range("S6").Value = [double]"0.0104386818071053"
$ws.Range("T6").Value = [double]"0.00675168658945601"

$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"11.717184"
$ws.Range("H7").Value = [double]"35.151552"
$ws.Range("I7").Value = [double]"0.3365934391853332"
$ws.Range("J7").Value = [double]"0.4228058096077277"
$ws.Range("K7").Value = [double]"2"
$ws.Range("M7").Value = [double]"35.789624"
$ws.Range("N7").Value = [double]"71.57924800000001"
$ws.Range("O7").Value = [double]"0.258139457682779"
$ws.Range("P7").Value = [double]"0.1993778771086309"
$ws.Range("Q7").Value = [double]"419.3536096988161"
$ws.Range("R7").Value = [double]"2516.121658192897"
$ws.Range("S7").Value = [double]"0.08688804785088336"
$ws.Range("T7").Value = [double]"0.08429812474878473"

$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"11.717184"
$ws.Range("H8").Value = [double]"35.151552"
$ws.Range("I8").Value = [double]"0.3365934391853332"
$ws.Range("J8").Value = [double]"0.4228058096077277"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"60.113367"
$ws.Range("N8").Value = [double]"180.340101"
$ws.Range("O8").Value = [double]"0.4335790718803266"
$ws.Range("P8").Value = [double]"0.5023219368682956"
$ws.Range("Q8").Value = [double]"704.3593819985281"
$ws.Range("R8").Value = [double]"6339.234437986753"
$ws.Range("S8").Value = [double]"0.1459398709629839"
$ws.Range("T8").Value = [double]"0.2123846332013216"

$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"11.717184"
$ws.Range("H9").Value = [double]"35.151552"
$ws.Range("I9").Value = [double]"0.3365934391853332"
$ws.Range("J9").Value = [double]"0.4228058096077277"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"0.2072186666666667"
$ws.Range("N9").Value = [double]"0.621656"
$ws.Range("O9").Value = [double]"0.001494603973349423"
$ws.Range("P9").Value = [double]"0.001731569652308208"
$ws.Range("Q9").Value = [double]"2.428019245568"
$ws.Range("R9").Value = [double]"21.852173210112"
$ws.Range("S9").Value = [double]"0.0005030738916097465"
$ws.Range("T9").Value = [double]"0.0007321177087363437"

$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"11.717184"
$ws.Range("H10").Value = [double]"35.151552"
$ws.Range("I10").Value = [double]"0.3365934391853332"
$ws.Range("J10").Value = [double]"0.4228058096077277"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"21.40334366666667"
$ws.Range("N10").Value = [double]"64.210031"
$ws.Range("O10").Value = [double]"0.1543756795743782"
$ws.Range("P10").Value = [double]"0.178851553034748"
$ws.Range("Q10").Value = [double]"250.7869159575681"
$ws.Range("R10").Value = [double]"2257.082243618112"
$ws.Range("S10").Value = [double]"0.05196184091451296"
$ws.Range("T10").Value = [double]"0.07561947568045607"

$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"11.717184"
$ws.Range("H11").Value = [double]"35.151552"
$ws.Range("I11").Value = [double]"0.3365934391853332"
$ws.Range("J11").Value = [double]"0.4228058096077277"
$ws.Range("K11").Value = [double]"2"
$ws.Range("M11").Value = [double]"21.1309775"
$ws.Range("N11").Value = [double]"42.261955"
$ws.Range("O11").Value = [double]"0.1524111868891667"
$ws.Range("P11").Value = [double]"0.1177170633360173"
$ws.Range("Q11").Value = [double]"247.59555146736"
$ws.Range("R11").Value = [double]"1485.57330880416"
$ws.Range("S11").Value = [double]"0.05130060556534318"
$ws.Range("T11").Value = [double]"0.04977145826842895"

$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"1.799402"
$ws.Range("H12").Value = [double]"5.398206"
$ws.Range("I12").Value = [double]"0.05169048362276865"
$ws.Range("J12").Value = [double]"0.06493007359274758"
$ws.Range("K12").Value = [double]"2"
$ws.Range("M12").Value = [double]"35.789624"
$ws.Range("N12").Value = [double]"71.57924800000001"
$ws.Range("O12").Value = [double]"0.258139457682779"
$ws.Range("P12").Value = [double]"0.1993778771086309"
$ws.Range("Q12").Value = [double]"64.39992100484801"
$ws.Range("R12").Value = [double]"386.3995260290881"
$ws.Range("S12").Value = [double]"0.01334335340974207"
$ws.Range("T12").Value = [double]"0.01294562023342919"

$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"1.799402"
$ws.Range("H13").Value = [double]"5.398206"
$ws.Range("I13").Value = [double]"0.05169048362276865"
$ws.Range("J13").Value = [double]"0.06493007359274758"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"60.113367"
$ws.Range("N13").Value = [double]"180.340101"
$ws.Range("O13").Value = [double]"0.4335790718803266"
$ws.Range("P13").Value = [double]"0.5023219368682956"
$ws.Range("Q13").Value = [double]"108.168112806534"
$ws.Range("R13").Value = [double]"973.513015258806"
$ws.Range("S13").Value = [double]"0.02241191191420526"
$ws.Range("T13").Value = [double]"0.03261580032810994"

$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"1.799402"
$ws.Range("H14").Value = [double]"5.398206"
$ws.Range("I14").Value = [double]"0.05169048362276865"
$ws.Range("J14").Value = [double]"0.06493007359274758"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"0.2072186666666667"
$ws.Range("N14").Value = [double]"0.621656"
$ws.Range("O14").Value = [double]"0.001494603973349423"
$ws.Range("P14").Value = [double]"0.001731569652308208"
$ws.Range("Q14").Value = [double]"0.3728696832373333"
$ws.Range("R14").Value = [double]"3.355827149136"
$ws.Range("S14").Value = [double]"7.725680220694331E-05"
$ws.Range("T14").Value = [double]"0.0001124309449553403"

$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"1.799402"
$ws.Range("H15").Value = [double]"5.398206"
$ws.Range("I15").Value = [double]"0.05169048362276865"
$ws.Range("J15").Value = [double]"0.06493007359274758"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"21.40334366666667"
$ws.Range("N15").Value = [double]"64.210031"
$ws.Range("O15").Value = [double]"0.1543756795743782"
$ws.Range("P15").Value = [double]"0.178851553034748"
$ws.Range("Q15").Value = [double]"38.51321940048734"
$ws.Range("R15").Value = [double]"346.618974604386"
$ws.Range("S15").Value = [double]"0.007979753536793177"
$ws.Range("T15").Value = [double]"0.01161284450072338"

$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"1.799402"
$ws.Range("H16").Value = [double]"5.398206"
$ws.Range("I16").Value = [double]"0.05169048362276865"
$ws.Range("J16").Value = [double]"0.06493007359274758"
$ws.Range("K16").Value = [double]"2"
$ws.Range("M16").Value = [double]"21.1309775"
$ws.Range("N16").Value = [double]"42.261955"
$ws.Range("O16").Value = [double]"0.1524111868891667"
$ws.Range("P16").Value = [double]"0.1177170633360173"
$ws.Range("Q16").Value = [double]"38.023123175455"
$ws.Range("R16").Value = [double]"228.13873905273"
$ws.Range("S16").Value = [double]"0.007878207959821202"
$ws.Range("T16").Value = [double]"0.00764337758552973"

$ws.Range("E17").Value = [double]"2"
$ws.Range("G17").Value = [double]"18.910282"
$ws.Range("H17").Value = [double]"37.820564"
$ws.Range("I17").Value = [double]"0.5432258172564757"
$ws.Range("J17").Value = [double]"0.4549089093375132"
$ws.Range("K17").Value = [double]"2"
$ws.Range("M17").Value = [double]"35.789624"
$ws.Range("N17").Value = [double]"71.57924800000001"
$ws.Range("O17").Value = [double]"0.258139457682779"
$ws.Range("P17").Value = [double]"0.1993778771086309"
$ws.Range("Q17").Value = [double]"676.7918825139681"
$ws.Range("R17").Value = [double]"2707.167530055872"
$ws.Range("S17").Value = [double]"0.140228017865871"
$ws.Range("T17").Value = [double]"0.09069877262151602"

$ws.Range("E18").Value = [double]"2"
$ws.Range("G18").Value = [double]"18.910282"
$ws.Range("H18").Value = [double]"37.820564"
$ws.Range("I18").Value = [double]"0.5432258172564757"
$ws.Range("J18").Value = [double]"0.4549089093375132"
$ws.Range("K18").Value = [double]"3"
$ws.Range("M18").Value = [double]"60.113367"
$ws.Range("N18").Value = [double]"180.340101"
$ws.Range("O18").Value = [double]"0.4335790718803266"
$ws.Range("P18").Value = [double]"0.5023219368682956"
$ws.Range("Q18").Value = [double]"1136.760721939494"
$ws.Range("R18").Value = [double]"6820.564331636963"
$ws.Range("S18").Value = [double]"0.2355313456674946"
$ws.Range("T18").Value = [double]"0.2285107244370635"

$ws.Range("E19").Value = [double]"2"
$ws.Range("G19").Value = [double]"18.910282"
$ws.Range("H19").Value = [double]"37.820564"
$ws.Range("I19").Value = [double]"0.5432258172564757"
$ws.Range("J19").Value = [double]"0.4549089093375132"
$ws.Range("K19").Value = [double]"3"
$ws.Range("M19").Value = [double]"0.2072186666666667"
$ws.Range("N19").Value = [double]"0.621656"
$ws.Range("O19").Value = [double]"0.001494603973349423"
$ws.Range("P19").Value = [double]"0.001731569652308208"
$ws.Range("Q19").Value = [double]"3.918563422330666"
$ws.Range("R19").Value = [double]"23.511380533984"
$ws.Range("S19").Value = [double]"0.0008119074648975162"
$ws.Range("T19").Value = [double]"0.0007877064619734639"

$ws.Range("E20").Value = [double]"2"
$ws.Range("G20").Value = [double]"18.910282"
$ws.Range("H20").Value = [double]"37.820564"
$ws.Range("I20").Value = [double]"0.5432258172564757"
$ws.Range("J20").Value = [double]"0.4549089093375132"
$ws.Range("K20").Value = [double]"3"
$ws.Range("M20").Value = [double]"21.40334366666667"
$ws.Range("N20").Value = [double]"64.210031"
$ws.Range("O20").Value = [double]"0.1543756795743782"
$ws.Range("P20").Value = [double]"0.178851553034748"
$ws.Range("Q20").Value = [double]"404.7432644795807"
$ws.Range("R20").Value = [double]"2428.459586877484"
$ws.Range("S20").Value = [double]"0.08386085470131542"
$ws.Range("T20").Value = [double]"0.08136116492435759"

$ws.Range("E21").Value = [double]"2"
$ws.Range("G21").Value = [double]"18.910282"
$ws.Range("H21").Value = [double]"37.820564"
$ws.Range("I21").Value = [double]"0.5432258172564757"
$ws.Range("J21").Value = [double]"0.4549089093375132"
$ws.Range("K21").Value = [double]"2"
$ws.Range("M21").Value = [double]"21.1309775"
$ws.Range("N21").Value = [double]"42.261955"
$ws.Range("O21").Value = [double]"0.1524111868891667"
$ws.Range("P21").Value = [double]"0.1177170633360173"
$ws.Range("Q21").Value = [double]"399.592743460655"
$ws.Range("R21").Value = [double]"1598.37097384262"
$ws.Range("S21").Value = [double]"0.082793691556897"
$ws.Range("T21").Value = [double]"0.05355054089260258"
